$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Add 6 new rows of annotations (rows 241-246) for the new lesson
# "97. Endpoints para buscar estados e cidades" under section 8.
# ------------------------------------------------------------------

$sectionName = "Ajustes finais no backend e no bucket"
$lessonName  = "Endpoints para buscar estados e cidades"

$notes = @(
  "2:07`n8. Ajustes finais no backend e bucket`n97. Endpoints para buscar estados e cidades`nutilização de padrão de nomes do spring data",
  "3:42`n8. Ajustes finais no backend e bucket`n97. Endpoints para buscar estados e cidades`nsempre importante lembrar de colocar a anotação @Service ... caso contrário os objetos não serão injetaveis",
  "3:52`n8. Ajustes finais no backend e bucket`n97. Endpoints para buscar estados e cidades`nIMPORTANTE: as classes services tem a unção de repassar a chamada para a camada repository",
  "4:08`n8. Ajustes finais no backend e bucket`n97. Endpoints para buscar estados e cidades`ncriação de endpoint EstadoResource para acessar os `"/estados`"",
  "9:13`n8. Ajustes finais no backend e bucket`n97. Endpoints para buscar estados e cidades`nsugestão do professor (implementado na aula e no projeto): criação do endpoint para obter cidades direto no EstadoResource invés de criar no resource de cidades",
  "10:05`n8. Ajustes finais no backend e bucket`n97. Endpoints para buscar estados e cidades`nIMPORTANTE: endpoint dentro de chaves {} significa que ele recebe um parametro de URL. Ex. `"/{estadoId}/cidades`" ... este caso recebe as cidades de acordo com o parametro estadoId, ou seja, recebe as cidades de acordo com o estado"
)

$heights = @(60, 75, 75, 60, 90, 105)

$lo = $ws.ListObjects.Item(1)

$newRows = @()

for ($i = 0; $i -lt $notes.Length; $i++) {
    $listRow = $lo.ListRows.Add()
    $r = $listRow.Range.Row
    $newRows += $r

    # Copy formatting from the row right above so styles (s="8","9","10")
    # and column formats match the rest of the table.
    $ws.Range("B" + ($r - 1) + ":G" + ($r - 1)).Copy()
    $ws.Range("B" + $r + ":G" + $r).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 2).Value = 8
    $ws.Cells.Item($r, 3).Value = $sectionName
    $ws.Cells.Item($r, 4).Value = 97
    $ws.Cells.Item($r, 6).Value = $notes[$i]

    $ws.Rows.Item($r).RowHeight = $heights[$i]
}

# Fill in the (repeated) lesson name last, across all the new rows, so the
# new shared string it introduces is appended after the six note strings.
foreach ($r in $newRows) {
    $ws.Cells.Item($r, 5).Value = $lessonName
}

$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Update the view state: new selection and scroll position
# ------------------------------------------------------------------
$ws.Range("F245").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 239
$win.ScrollColumn = 1
